# xHI_EoR.xlsx update — reclassify "Lya emitting galaxies" method rows into
# more specific categories, and append three new rows of data from Tang+24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reclassify existing "Lya emitting galaxies" Method entries (column I) ---
# Mason+18 / Mason+19 -> Lya break galaxies
$ws.Cells.Item(9, 9).Value  = "Lya break galaxies"
$ws.Cells.Item(11, 9).Value = "Lya break galaxies"

# Goto+21 / Morales+21 / Wold+22 -> Lya LF
$ws.Cells.Item(16, 9).Value = "Lya LF"
$ws.Cells.Item(17, 9).Value = "Lya LF"
$ws.Cells.Item(18, 9).Value = "Lya LF"
$ws.Cells.Item(19, 9).Value = "Lya LF"
$ws.Cells.Item(24, 9).Value = "Lya LF"

# Bruton+23 -> Lya EW
$ws.Cells.Item(28, 9).Value = "Lya EW"

# --- Clean up stray formatting left over on rows 2 and 20:22 (Normal style) ---
$ws.Range("A2:J2").Style = "Normal"
$ws.Range("A20:J22").Style = "Normal"
$ws.Range("E20:F22").NumberFormat = "0.000"

# --- Append three new rows (Tang+24) ---
# Row 54: z = 10.0-13.3
$ws.Cells.Item(54, 1).Value = "10.0-13.3"
$ws.Cells.Item(54, 4).Value = 0.89
$ws.Cells.Item(54, 5).Value = 0.08
$ws.Cells.Item(54, 6).Value = 0.21
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = "Lya EW"
$ws.Cells.Item(54, 10).Value = "Tang+24"

# Row 55: z = 8.0-10.0
$ws.Cells.Item(55, 1).Value = "8.0-10.0"
$ws.Cells.Item(55, 4).Value = 0.81
$ws.Cells.Item(55, 5).Value = 0.12
$ws.Cells.Item(55, 6).Value = 0.24
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = "Lya EW"
$ws.Cells.Item(55, 10).Value = "Tang+24"

# Row 56: z = 6.5-8.0
$ws.Cells.Item(56, 1).Value = "6.5-8.0"
$ws.Cells.Item(56, 4).Value = 0.48
$ws.Cells.Item(56, 5).Value = 0.15
$ws.Cells.Item(56, 6).Value = 0.22
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = "Lya EW"
$ws.Cells.Item(56, 10).Value = "Tang+24"

# A55 picked up a (cosmetic) date-style number format in the source edit even
# though it holds text.
$ws.Range("A55").NumberFormat = "d-mmm"

# --- Update view / selection state to match the edited workbook ---
$ws.Range("F57").Select()
